$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 242.77777
$ws.Range("I4").Value = 254.375
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 254.375
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -140.375
$ws.Range("N4").Value = -378
$ws.Range("H86").Value = 1553.75
$ws.Range("J86").Value = 1808.3334
$ws.Range("L86").Value = 1808.3334
$ws.Range("N86").Value = -4054.3334
$ws.Range("H89").Value = 1553.75
$ws.Range("J89").Value = 1808.3334
$ws.Range("L89").Value = 9041.666999999999
$ws.Range("N89").Value = -20273.667
$ws.Range("H98").Value = 3414.5
$ws.Range("I98").Value = 3189
$ws.Range("K98").Value = 3189
$ws.Range("M98").Value = -1691
$ws.Range("H107").Value = 72.5
$ws.Range("I107").Value = 72.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 72.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1847.5
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 1750
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1754
$ws.Range("N113").Value = -8508
$ws.Range("H122").Value = 3414.5
$ws.Range("I122").Value = 3189
$ws.Range("K122").Value = 9567
$ws.Range("M122").Value = -7117
$ws.Range("H132").Value = 5132.2104
$ws.Range("I132").Value = 4912.4707
$ws.Range("K132").Value = 14737.4121
$ws.Range("M132").Value = -12207.4121
$ws.Range("H135").Value = 3998.6667
$ws.Range("I135").Value = 3999
$ws.Range("J135").Value = 3998
$ws.Range("K135").Value = 35991
$ws.Range("L135").Value = 35982
$ws.Range("M135").Value = -33456
$ws.Range("N135").Value = -41052
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H140").Value = 70001
$ws.Range("J140").Value = 70001
$ws.Range("L140").Value = 70001
$ws.Range("N140").Value = -80361

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6863.75
$ws.Range("I32").Value = 6863.75
$ws.Range("K32").Value = 6863.75
$ws.Range("M32").Value = -6576.75
$ws.Range("H45").Value = 2685.5
$ws.Range("I45").Value = 2558.4
$ws.Range("J45").Value = 2897.3333
$ws.Range("K45").Value = 2558.4
$ws.Range("L45").Value = 2897.3333
$ws.Range("M45").Value = -2181.4
$ws.Range("N45").Value = -3651.3333
$ws.Range("H88").Value = 619
$ws.Range("J88").Value = 665
$ws.Range("L88").Value = 665
$ws.Range("N88").Value = -1477
$ws.Range("H91").Value = 619
$ws.Range("J91").Value = 665
$ws.Range("L91").Value = 665
$ws.Range("N91").Value = -3473
$ws.Range("H132").Value = 3541.2
$ws.Range("I132").Value = 3541.2
$ws.Range("K132").Value = 10623.6
$ws.Range("M132").Value = -8093.599999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2611.5
$ws.Range("I20").Value = 1769.25
$ws.Range("J20").Value = 3453.75
$ws.Range("K20").Value = 1769.25
$ws.Range("L20").Value = 3453.75
$ws.Range("M20").Value = -1522.25
$ws.Range("N20").Value = -3947.75
$ws.Range("H105").Value = 76500
$ws.Range("I105").Value = 150000
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 150000
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -148253
$ws.Range("N105").Value = -6494
$ws.Range("H134").Value = 2125.8235
$ws.Range("I134").Value = 2139.8125
$ws.Range("K134").Value = 6419.4375
$ws.Range("M134").Value = -3884.4375

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 80.55556
$ws.Range("J7").Value = 54
$ws.Range("L7").Value = 54
$ws.Range("N7").Value = -280
$ws.Range("H62").Value = 8202.4
$ws.Range("J62").Value = 9001
$ws.Range("L62").Value = 9001
$ws.Range("N62").Value = -10249
$ws.Range("H65").Value = 8202.4
$ws.Range("J65").Value = 9001
$ws.Range("L65").Value = 45005
$ws.Range("N65").Value = -51245
$ws.Range("H120").Value = 49999
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 49999
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 49999
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -57257

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 834.3
$ws.Range("I4").Value = 49.5
$ws.Range("K4").Value = 148.5
$ws.Range("M4").Value = -36.5
$ws.Range("H97").Value = 23330
$ws.Range("I97").Value = 9995
$ws.Range("J97").Value = 50000
$ws.Range("K97").Value = 29985
$ws.Range("L97").Value = 150000
$ws.Range("M97").Value = -29489
$ws.Range("N97").Value = -150992
$ws.Range("H113").Value = 536.2222
$ws.Range("I113").Value = 540.6667
$ws.Range("J113").Value = 534
$ws.Range("K113").Value = 1622.0001
$ws.Range("L113").Value = 1602
$ws.Range("M113").Value = 547.9999
$ws.Range("N113").Value = -5942
$ws.Range("H138").Value = 2000
$ws.Range("J138").Value = 2000
$ws.Range("L138").Value = 6000
$ws.Range("N138").Value = -16280
$ws.Range("H140").Value = 2337
$ws.Range("I140").Value = 2337
$ws.Range("K140").Value = 7011
$ws.Range("M140").Value = -1831

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 2111
$ws.Range("I122").Value = 1355.6
$ws.Range("K122").Value = 4066.8
$ws.Range("M122").Value = -1616.8
$ws.Range("H132").Value = 4631.75
$ws.Range("I132").Value = 4631.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13895.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11365.25
$ws.Range("N132").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 17218.666
$ws.Range("I136").Value = 3351.5
$ws.Range("K136").Value = 10054.5
$ws.Range("M136").Value = -7504.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 300
$ws.Range("I107").Value = 300
$ws.Range("K107").Value = 900
$ws.Range("M107").Value = 1020
$ws.Range("H122").Value = 1393
$ws.Range("I122").Value = 1393
$ws.Range("K122").Value = 4179
$ws.Range("M122").Value = -1729
$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H132").Value = 1996.4
$ws.Range("I132").Value = 1497
$ws.Range("K132").Value = 4491
$ws.Range("M132").Value = -1961
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360
